# Generate Report for Handback
# Adds "Latest Target File" (F) and "Latest Handback File" (G) hyperlinked
# entries for rows 2-3 on the zh-cn and de-de sheets, updates the status
# text to reflect a completed handback, and refreshes the handback
# datetime stamps.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# --- Update status text everywhere it currently reads "Ready for handoff" ---
$ws1.Range("B2").Value = $statusText
$ws1.Range("C2").Value = $statusText
$ws1.Range("B3").Value = $statusText
$ws1.Range("C3").Value = $statusText

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

# --- zh-cn sheet: add Latest Target File (F) / Latest Handback File (G) ---
$zhAUrl = "https://github.com/OpenLocalizationTest/oltest/blob/9005276061917fd52f340e21a110314c9adf4b38/e2e/a.md"
$zhBUrl = "https://github.com/OpenLocalizationTest/oltest/blob/9005276061917fd52f340e21a110314c9adf4b38/e2e/b.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e2e6ee3b685cd17ee844665ce942208f544d05c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

# Rebuild the hyperlink collection so the new Target/Handback links land in
# document order right after the existing row's Handoff File link (matching
# how Excel lays out newly-inserted columns within a row).
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhAUrl, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhXlfUrl, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhAUrl, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhXlfUrl, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $zhBUrl, "", "", "b.md")
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhXlfUrl, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhAUrl, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhXlfUrl, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")

# Keep every hyperlinked cell on this sheet using the workbook's hyperlink
# font color (matches the pre-existing Source/Handoff File link styling).
$wsZh.Range("A2:A3").Font.Color = 15570276
$wsZh.Range("D2:D3").Font.Color = 15570276
$wsZh.Range("F2:G3").Font.Color = 15570276

# Latest Handback DateTime for zh-cn
$wsZh.Range("H2").Value = "2016-03-22 06:35:10"
$wsZh.Range("H3").Value = "2016-03-22 06:35:10"

# --- de-de sheet: add Latest Target File (F) / Latest Handback File (G) ---
$deAUrl = "https://github.com/OpenLocalizationTest/oltest/blob/9005276061917fd52f340e21a110314c9adf4b38/e2e/a.md"
$deBUrl = "https://github.com/OpenLocalizationTest/oltest/blob/9005276061917fd52f340e21a110314c9adf4b38/e2e/b.md"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/93c21e2d0b38a0d3cfd7c12b00f10b5c04781874/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deAUrl, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deXlfUrl, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deAUrl, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deXlfUrl, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $deBUrl, "", "", "b.md")
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deXlfUrl, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deAUrl, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deXlfUrl, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")

$wsDe.Range("A2:A3").Font.Color = 15570276
$wsDe.Range("D2:D3").Font.Color = 15570276
$wsDe.Range("F2:G3").Font.Color = 15570276

# Latest Handback DateTime for de-de (distinct stamp from zh-cn)
$wsDe.Range("H2").Value = "2016-03-22 06:35:17"
$wsDe.Range("H3").Value = "2016-03-22 06:35:17"
